$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4602
$ws.Range("K3").Value = 4721
$ws.Range("K4").Value = 960
$ws.Range("K5").Value = 343
$ws.Range("J6").Value = 11054
$ws.Range("K6").Value = 5328
$ws.Range("J7").Value = 29300
$ws.Range("K7").Value = 15954

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 140
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 473
$ws.Range("K8").Value = 1070
$ws.Range("K12").Value = 28
$ws.Range("K15").Value = 163
$ws.Range("K18").Value = 109
$ws.Range("K19").Value = 479
$ws.Range("K20").Value = 369
$ws.Range("K24").Value = 47
$ws.Range("K29").Value = 847
$ws.Range("K30").Value = 59
$ws.Range("J31").Value = 305
$ws.Range("K31").Value = 182
$ws.Range("K33").Value = 669
$ws.Range("K36").Value = 202
$ws.Range("K37").Value = 537
$ws.Range("K41").Value = 122
$ws.Range("K42").Value = 593
$ws.Range("K51").Value = 201
$ws.Range("K52").Value = 421
$ws.Range("K54").Value = 302
$ws.Range("K55").Value = 184
$ws.Range("K57").Value = 54
$ws.Range("K65").Value = 366
$ws.Range("K67").Value = 615
$ws.Range("K68").Value = 41
$ws.Range("K77").Value = 116
$ws.Range("K79").Value = 396
$ws.Range("K83").Value = 340
$ws.Range("K84").Value = 117
$ws.Range("K85").Value = 726
$ws.Range("K86").Value = 105
$ws.Range("K89").Value = 228
$ws.Range("K90").Value = 145
$ws.Range("K91").Value = 172
$ws.Range("K94").Value = 201
$ws.Range("K95").Value = 280
$ws.Range("K99").Value = 268
$ws.Range("J101").Value = 29300
$ws.Range("K101").Value = 15954

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 166
$ws.Range("K3").Value = 151
$ws.Range("K7").Value = 473

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 253
$ws.Range("K3").Value = 243
$ws.Range("K6").Value = 167
$ws.Range("K7").Value = 726

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 113
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 421

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 323
$ws.Range("K7").Value = 1070

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 121
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 183
$ws.Range("K6").Value = 192
$ws.Range("K7").Value = 669

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 95
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 151
$ws.Range("K6").Value = 160
$ws.Range("K7").Value = 537

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 112
$ws.Range("K7").Value = 366

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K5").Value = 9
$ws.Range("K7").Value = 268

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 110
$ws.Range("K6").Value = 69
$ws.Range("J7").Value = 305
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 179
$ws.Range("K6").Value = 179
$ws.Range("K7").Value = 615

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 50
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 153
$ws.Range("K7").Value = 302

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 243
$ws.Range("K3").Value = 300
$ws.Range("K5").Value = 25
$ws.Range("K7").Value = 847

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 150
$ws.Range("K7").Value = 479

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 223
$ws.Range("K7").Value = 593

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K4").Value = 26
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 396

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 369

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 58
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 42
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 28
